$d = $word.ActiveDocument

# 1. Split the leading space off the Title run and insert an empty
#    "_GoBack" bookmark right after it (mirrors Word's automatic
#    last-edit-location bookmark). Character 1 is immediately after the
#    leading space of " Title ".
$goBackRange = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# 2. Force the two "DefinitionTerm" runs to be rebuilt so the stale
#    lastRenderedPageBreak hint Word had cached is dropped.
$d.Content.Find.Execute(" DefinitionTerm ", $true, $false, $false, $false, $false, $true, 1, $false, " DefinitionTerm ", 2) | Out-Null

# 3. Shrink the two-column table slightly (895 -> 842 twips, i.e.
#    44.75pt -> 42.1pt).
$t = $d.Tables(1)
$t.Columns(1).Width = 42.1
$t.Columns(2).Width = 42.1

# 4. Reduce the "Body Text" (Corpotesto) and "Block Text" (Testodelblocco)
#    styles' character size from 13pt/12pt down to the document default
#    of 12pt.
$d.Styles("Corpotesto").Font.Size = 12
$d.Styles("Testodelblocco").Font.Size = 12
